# Sync attendance_reports: reorder "Recorded By" (column G) values so that
# when the list starts with "System, ", the "System" entry is moved to the
# end of the comma-separated list (e.g. "System, foo@bar.com" becomes
# "foo@bar.com, System").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 157 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) { continue }

    if ($text -like "System, *") {
        $parts = $text -split ", "
        if ($parts.Count -gt 1) {
            $rest = $parts[1..($parts.Count - 1)]
            $newText = ($rest -join ", ") + ", System"
            $cell.Value = $newText
        }
    }
}
